{"js": "const replacements = [\n  [\"100\u00f79=\", \"137\u00f75=\"],\n  [\"222\u00f78=\", \"818\u00f76=\"],\n  [\"554\u00f74=\", \"432\u00f72=\"],\n  [\"147\u00f74=\", \"188\u00f74=\"],\n  [\"467\u00f75=\", \"972\u00f79=\"],\n  [\"519\u00f72=\", \"398\u00f72=\"],\n  [\"102\u00f78=\", \"671\u00f72=\"],\n  [\"450\u00f72=\", \"550\u00f72=\"],\n  [\"431\u00f76=\", \"639\u00f75=\"],\n  [\"774\u00f76=\", \"830\u00f75=\"],\n  [\"271\u00f74=\", \"725\u00f73=\"],\n  [\"918\u00f74=\", \"712\u00f77=\"],\n  [\"514\u00f79=\", \"330\u00f78=\"],\n  [\"285\u00f73=\", \"106\u00f78=\"],\n  [\"880\u00f78=\", \"452\u00f78=\"],\n  [\"756\u00f76=\", \"693\u00f77=\"],\n  [\"843\u00f73=\", \"585\u00f79=\"],\n  [\"610\u00f74=\", \"861\u00f76=\"],\n  [\"853\u00f77=\", \"992\u00f74=\"],\n  [\"802\u00f76=\", \"817\u00f76=\"],\n  [\"484\u00f78=\", \"173\u00f76=\"],\n  [\"392\u00f74=\", \"958\u00f78=\"],\n  [\"182\u00f74=\", \"633\u00f74=\"],\n  [\"752\u00f74=\", \"568\u00f75=\"],\n  [\"775\u00f78=\", \"904\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"100\u00f79=\", \"137\u00f75=\"),\n    @(\"222\u00f78=\", \"818\u00f76=\"),\n    @(\"554\u00f74=\", \"432\u00f72=\"),\n    @(\"147\u00f74=\", \"188\u00f74=\"),\n    @(\"467\u00f75=\", \"972\u00f79=\"),\n    @(\"519\u00f72=\", \"398\u00f72=\"),\n    @(\"102\u00f78=\", \"671\u00f72=\"),\n    @(\"450\u00f72=\", \"550\u00f72=\"),\n    @(\"431\u00f76=\", \"639\u00f75=\"),\n    @(\"774\u00f76=\", \"830\u00f75=\"),\n    @(\"271\u00f74=\", \"725\u00f73=\"),\n    @(\"918\u00f74=\", \"712\u00f77=\"),\n    @(\"514\u00f79=\", \"330\u00f78=\"),\n    @(\"285\u00f73=\", \"106\u00f78=\"),\n    @(\"880\u00f78=\", \"452\u00f78=\"),\n    @(\"756\u00f76=\", \"693\u00f77=\"),\n    @(\"843\u00f73=\", \"585\u00f79=\"),\n    @(\"610\u00f74=\", \"861\u00f76=\"),\n    @(\"853\u00f77=\", \"992\u00f74=\"),\n    @(\"802\u00f76=\", \"817\u00f76=\"),\n    @(\"484\u00f78=\", \"173\u00f76=\"),\n    @(\"392\u00f74=\", \"958\u00f78=\"),\n    @(\"182\u00f74=\", \"633\u00f74=\"),\n    @(\"752\u00f74=\", \"568\u00f75=\"),\n    @(\"775\u00f78=\", \"904\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
